# ProjectFundingData.xlsx edit: create the "regular proposal" page (Sheet2)
# by copying Sheet1's funding table and appending a "SpeedUp" column,
# then making Sheet2 the active tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Column widths on Sheet2 (same best-fit widths as Sheet1 cols A:H) ---
$ws2.Columns.Item(1).ColumnWidth = 16
$ws2.Columns.Item(2).ColumnWidth = 12.42578125
$ws2.Columns.Item(3).ColumnWidth = 16.42578125
$ws2.Columns.Item(4).ColumnWidth = 18.7109375
$ws2.Columns.Item(5).ColumnWidth = 19.5703125
$ws2.Columns.Item(6).ColumnWidth = 15
$ws2.Columns.Item(7).ColumnWidth = 16
$ws2.Columns.Item(8).ColumnWidth = 12.42578125

# --- Header row (row 1): same labels as Sheet1, plus new "SpeedUp" header in I1 ---
$ws2.Range("A1").Value = "ExistingProjFund"
$ws2.Range("B1").Value = "ReqProjFund"
$ws2.Range("C1").Value = "BudgetPersonnel"
$ws2.Range("D1").Value = "BudgetPatientCosts"
$ws2.Range("E1").Value = "BudgetCoreFacilities"
$ws2.Range("F1").Value = "BudgetSupplies"
$ws2.Range("G1").Value = "BudgetOmicData"
$ws2.Range("H1").Value = "BudgetOther"
$ws2.Range("I1").Value = "SpeedUp"

# --- Data row (row 2): same values as Sheet1, plus the new description text in I2 ---
$ws2.Range("A2").Value = 1000
$ws2.Range("B2").Value = 2000
$ws2.Range("C2").Value = "119994"
$ws2.Range("D2").Value = "6000"
$ws2.Range("E2").Value = 0
$ws2.Range("F2").Value = "47250"
$ws2.Range("G2").Value = 0
$ws2.Range("H2").Value = "126756"
$ws2.Range("I2").Value = "`"Lorem ipsum dolor sit amet, consectetur adipiscing elit, sed do eiusmod tempor incididunt ut labore et dolore magna aliqua. Ut enim ad minim veniam, quis nostrud exercitation ullamco laboris nisi ut aliquip ex ea commodo consequat. Duis aute irure dolor in reprehenderit in voluptate velit esse cillum dolore eu fugiat nulla pariatur. Excepteur sint occaecat cupidatat non proident, sunt in culpa qui officia deserunt mollit anim id est laborum.`""

# New cell (I2) uses a distinct font: Arial, black, size 11
$ws2.Range("I2").Font.Name = "Arial"
$ws2.Range("I2").Font.Color = 0

# Page setup: portrait orientation on Sheet2
$ws2.PageSetup.Orientation = 1

# --- Selection / view state ---
$ws1.Range("A1:H2").Select()
$ws2.Range("I2").Select()
$ws2.Activate()

Write-Host "edit complete"
